# Added Veenkampen BAM data (18 new streams), added BAM_Conc to Veenkampen
# Air Quality package.
#
# This script appends 18 new instrument rows (BAM_* streams) to the
# "Veenkampen" worksheet (rows 195-212) and records the edit in the
# "Metadata" worksheet's changelog (row 20), matching two earlier log
# entries' dates being corrected as well.

$wb = $excel.ActiveWorkbook
$wsVeenkampen = $wb.Worksheets.Item("Veenkampen")
$wsMetadata   = $wb.Worksheets.Item("Metadata")

# ---------------------------------------------------------------------------
# 1. New instrument rows on the Veenkampen sheet (rows 195-212)
#    Columns: A=Stream name, B=Unit, C=Instrument, D=Long name,
#             F=Stream, G=Measurement interval
# ---------------------------------------------------------------------------

$instrument = "Beta Attenuation Mass Monitor 1020"
$interval   = "1 hour"

$rows = @(
    @{ Row=195; A="BAM_Conc";  B="mg m-3"; Dfont4=$false; D="Particulate matter concentration < 2.5 µm";                 F="SINGLE/ADVANCED-AIR QUALITY" },
    @{ Row=196; A="BAM_Qtot";  B="m3";     Dfont4=$true;  D="Total flow volume";                                          F="SINGLE" },
    @{ Row=197; A="BAM_WS";    B="m s-1";  Dfont4=$false; D="Wind Speed (BX-591)";                                        F="SINGLE" },
    @{ Row=198; A="BAM_RH";    B="%";      Dfont4=$true;  D="Relative humidity (BX-593)";                                 F="SINGLE" },
    @{ Row=199; A="BAM_Delta"; B="°C";     Dfont4=$true;  D="Delta temperature (BX-597)";                                 F="SINGLE" },
    @{ Row=200; A="BAM_AT";    B="°C";     Dfont4=$true;  D="Air temperature (BX-597)";                                   F="SINGLE" },
    @{ Row=201; A="BAM_E";     B="-";      Dfont4=$false; D="Flag: External reset or Interface Reset";                    F="SINGLE" },
    @{ Row=202; A="BAM_U";     B="-";      Dfont4=$true;  D="Flag: Telemetry Fault or Interface Fault ";                  F="SINGLE" },
    @{ Row=203; A="BAM_M";     B="-";      Dfont4=$true;  D="Flag: Maintenance Alarm";                                    F="SINGLE" },
    @{ Row=204; A="BAM_I";     B="-";      Dfont4=$true;  D="Flag: Internal Error or Coarse Link Down";                   F="SINGLE" },
    @{ Row=205; A="BAM_L";     B="-";      Dfont4=$true;  D="Flag: Power Failure or Processor Reset ";                   F="SINGLE" },
    @{ Row=206; A="BAM_R";     B="-";      Dfont4=$true;  D="Flag: Reference Error or Membrane Timeout ";                 F="SINGLE" },
    @{ Row=207; A="BAM_N";     B="-";      Dfont4=$true;  D="Flag: Nozzle Error";                                         F="SINGLE" },
    @{ Row=208; A="BAM_F";     B="-";      Dfont4=$true;  D="Flag: Flow Error";                                           F="SINGLE" },
    @{ Row=209; A="BAM_P";     B="-";      Dfont4=$true;  D="Flag: Pressure Drop Alarm or Delta-Pressure Alarm ";         F="SINGLE" },
    @{ Row=210; A="BAM_D";     B="-";      Dfont4=$true;  D="Flag: Deviant Membrane Density Alarm or BAM CAL alarm";      F="SINGLE" },
    @{ Row=211; A="BAM_C";     B="-";      Dfont4=$true;  D="Flag: Count Error or Data Error ";                          F="SINGLE" },
    @{ Row=212; A="BAM_T";     B="-";      Dfont4=$true;  D="Flag: Tape System Error or Filter Tape Error ";              F="SINGLE" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # A: stream name
    $wsVeenkampen.Cells.Item($rowNum, 1).Value2 = $r.A

    # B: unit (plain string here; superscript units are fixed up below)
    $cellB = $wsVeenkampen.Cells.Item($rowNum, 2)
    $cellB.Value2 = $r.B
    $cellB.Font.Name = "Calibri"
    $cellB.Font.Size = 11

    # C: instrument
    $cellC = $wsVeenkampen.Cells.Item($rowNum, 3)
    $cellC.Value2 = $instrument
    $cellC.Font.Name = "Calibri"
    $cellC.Font.Size = 11

    # D: long name
    $cellD = $wsVeenkampen.Cells.Item($rowNum, 4)
    $cellD.Value2 = $r.D
    if ($r.Dfont4) {
        $cellD.Font.Name = "Calibri"
        $cellD.Font.Size = 11
    }

    # F: stream
    $wsVeenkampen.Cells.Item($rowNum, 6).Value2 = $r.F

    # G: measurement interval
    $wsVeenkampen.Cells.Item($rowNum, 7).Value2 = $interval
}

# Fix up the two units that need a superscripted exponent:
# B195 -> "mg m-3" with the "-3" superscripted
# B196 -> "m3" with the "3" superscripted
$charsConc = $wsVeenkampen.Range("B195").Characters(5, 2)
$charsConc.Font.Superscript = $true

$charsQtot = $wsVeenkampen.Range("B196").Characters(2, 1)
$charsQtot.Font.Superscript = $true

# ---------------------------------------------------------------------------
# 2. Metadata changelog updates
# ---------------------------------------------------------------------------

# Correct the dates of the three most recent pre-existing log entries
# (rows 13-15) back by 366 days.
$wsMetadata.Cells.Item(13, 2).Value2 = 45250
$wsMetadata.Cells.Item(14, 2).Value2 = 45261
$wsMetadata.Cells.Item(15, 2).Value2 = 45268

# Add the new changelog entry for this edit in row 20.
$wsMetadata.Cells.Item(20, 2).Value2 = 45638
$wsMetadata.Cells.Item(20, 2).NumberFormat = "m/d/yy"
$wsMetadata.Cells.Item(20, 3).Value2 = "Sjoerd Barten"
$wsMetadata.Cells.Item(20, 4).Value2 = "Added Veenkampen BAM data (18 new streams), added BAM_Conc to Veenkampen Air Quality package"
